# Single function output now working with named areas
#
# Adds a 4-cell "inputs" row (C1:F1) and a 4-cell "outputs" row (C2:F2) to
# the existing single-input/single-output example sheet, wires them up with
# formulas, and registers the two new ranges as workbook-level named areas
# (alongside the existing "input"/"output" names) so a caller can address
# the whole input/output vector by name.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New sample input values across C1:F1.
$ws.Range("C1").Value = 10
$ws.Range("D1").Value = 20
$ws.Range("E1").Value = 30
$ws.Range("F1").Value = 40

# D2:F2 share the same formula pattern as the original B2 (EXP(-1/x)),
# entered as one multi-cell assignment so it is stored as a shared formula.
$ws.Range("D2:F2").Formula = "=EXP(-1/D1)"

# C2 derives from the last of the shared-formula outputs.
$ws.Range("C2").Formula = "=F2-10"

# Register the new input/output vectors as workbook-level named ranges,
# matching the existing "input" / "output" single-cell names.
$wb.Names.Add("inputs", "=Sheet1!`$C`$1:`$F`$1")
$wb.Names.Add("outputs", "=Sheet1!`$C`$2:`$F`$2")

# Leave the selection parked below the new input block.
$ws.Range("C3").Select() | Out-Null

# Page setup tweaks that came along with this edit.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
